# Generate Report for Handoff
# Adds two newly-handed-off files (3d910c7f... and bd85d307...) to the
# Overview / zh-cn / de-de sheets, growing each table from 3 to 5 rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (table3 / "Overview")
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

$wsOverview.Range("A4").Value = "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-28 02:40:34"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md", "", "", "e2e\3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md")

$wsOverview.Range("A5").Value = "bd85d307-11d0-4207-8838-4c6012e1889a.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-28 02:40:34"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/bd85d307-11d0-4207-8838-4c6012e1889a.md", "", "", "e2e\bd85d307-11d0-4207-8838-4c6012e1889a.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table1 / "zh_cn")
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P5"))

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'False"
$wsZhCn.Range("G4").Value = "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.d62eaa93f590c0d247a28b32abb24731821c9c8b.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-28 02:40:29"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md", "", "", "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md")

$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "'False"
$wsZhCn.Range("G5").Value = "bd85d307-11d0-4207-8838-4c6012e1889a.ded0bf512234fe44a8e2b2ec0e81482c79ea91bb.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-28 02:40:29"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M5").Value = "'True"
$wsZhCn.Range("O5").Value = "'False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/bd85d307-11d0-4207-8838-4c6012e1889a.md", "", "", "bd85d307-11d0-4207-8838-4c6012e1889a.md")

# ---------------------------------------------------------------------
# Sheet "de-de" (table2 / "de_de")
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P5"))

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'False"
$wsDeDe.Range("G4").Value = "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.d62eaa93f590c0d247a28b32abb24731821c9c8b.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-28 02:40:34"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md", "", "", "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md")

$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "'False"
$wsDeDe.Range("G5").Value = "bd85d307-11d0-4207-8838-4c6012e1889a.ded0bf512234fe44a8e2b2ec0e81482c79ea91bb.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-08-28 02:40:34"
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M5").Value = "'True"
$wsDeDe.Range("O5").Value = "'False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/bd85d307-11d0-4207-8838-4c6012e1889a.md", "", "", "bd85d307-11d0-4207-8838-4c6012e1889a.md")

Write-Host "Handoff report rows added"
